# Update the error-table values for Q0..Q9 (rows 2-11) on the sheet.
# Column A labels (Q0..Q9) are unchanged; columns B-G (ME, MAE, MSE, RMSE, SE, N)
# are refreshed to reflect the full ifoCAST series evaluation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1979813930710988, 0.3795265156019396, 0.2778640158610643, 0.5271280829751572, 0.5069777247436844, 14)
    3  = @(0.2347256923014825, 0.3822201530932461, 0.2247519983339054, 0.4740801602407608, 0.428711832139554,  13)
    4  = @(0.2681184607162052, 0.3968952304966916, 0.2554487726035262, 0.5054194026781385, 0.4474915096726709, 12)
    5  = @(0.3296070495036373, 0.3883960298041425, 0.2818320098626162, 0.5308785264659103, 0.4364748825054373, 11)
    6  = @(0.3168892358128843, 0.3761758563838308, 0.2060985664200785, 0.4539807996161054, 0.342668901818269,  10)
    7  = @(0.3050354114201801, 0.3134960335622032, 0.2324223121793234, 0.4821019728017334, 0.3959768600611978,  9)
    8  = @(0.3287758598977879, 0.4033255329523637, 0.2534286472417211, 0.5034169715471669, 0.4075502860334328,  8)
    9  = @(0.3037301248181886, 0.4152369558703082, 0.3412314012878718, 0.5841501530324817, 0.5389582370896339,  7)
    10 = @(0.2667014930737632, 0.3376616306061401, 0.1853174931998413, 0.4304851834846831, 0.3701693776509343,  6)
    11 = @(0.3157875881341893, 0.415763948321914,  0.4482400691008893, 0.6695073331195778, 0.6600362379078902,  5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
